$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the third worksheet ("Module-name" -> "AddingResturant")
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "AddingResturant"

# ---------------------------------------------------------------------------
# Helper functions: write a value into a cell and then stamp it with the
# formatting (border/number-format/font) used elsewhere in the workbook by
# copying the format only (xlPasteFormats) from a donor cell that already
# carries the desired style.
# ---------------------------------------------------------------------------
function Set-PlainCell($cellRef, $value) {
    $ws3.Range($cellRef).Value = $value
    $ws1.Range("A1").Copy()
    $ws3.Range($cellRef).PasteSpecial(-4122)
}

function Set-HyperlinkStyledCell($cellRef, $value) {
    $ws3.Range($cellRef).Value = $value
    $ws1.Range("B2").Copy()
    $ws3.Range($cellRef).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 2. Populate the new "AddingResturant" sheet with the POM admin-page data
# ---------------------------------------------------------------------------
Set-PlainCell "A1" "UN"
Set-PlainCell "B1" "admin"

Set-PlainCell "A2" "PWD"
Set-PlainCell "B2" "codeastro"

Set-PlainCell "A3" "URL"
Set-PlainCell "B3" "http://rmgtestingserver/domain/Online_Food_Ordering_System/admin"

Set-PlainCell "A4" "Dropdownvalue"
Set-PlainCell "B4" "10am"

Set-PlainCell "A5" "Dropdowntext"
Set-PlainCell "B5" "Mon-Fri"

Set-PlainCell "A6" "Dropdownvalue2"
Set-PlainCell "B6" "Indian "

Set-PlainCell "A7" "Dropdowntext2"
Set-PlainCell "B7" "6pm"

Set-PlainCell "A8" "Address"
Set-PlainCell "B8" "Bangalore"

Set-PlainCell "A9" "Resturant_name"
Set-PlainCell "B9" "Navarang resturant"

Set-PlainCell "A10" "mob_no"
$ws3.Range("B10").Value = 9902990508
$ws1.Range("A1").Copy()
$ws3.Range("B10").PasteSpecial(-4122)

Set-PlainCell "A11" "Email"
Set-PlainCell "B11" "manuvirat775@gmail.com"

Set-PlainCell "A12" "Website"
Set-PlainCell "B12" "www.https.com"

# ---------------------------------------------------------------------------
# 3. Hyperlinks.
#    B11 (mail) and B12 (web) keep the plain "s=2" look, B3 (admin url) gets
#    the blue/underlined "Hyperlink" look ("s=4"), matching the target file.
# ---------------------------------------------------------------------------
$ws3.Hyperlinks.Add($ws3.Range("B11"), "mailto:manuvirat775@gmail.com")
Set-PlainCell "B11" "manuvirat775@gmail.com"

$ws3.Hyperlinks.Add($ws3.Range("B12"), "http://www.https.com")
Set-PlainCell "B12" "www.https.com"

$ws3.Hyperlinks.Add($ws3.Range("B3"), "http://rmgtestingserver/domain/Online_Food_Ordering_System/admin")
Set-HyperlinkStyledCell "B3" "http://rmgtestingserver/domain/Online_Food_Ordering_System/admin"

# ---------------------------------------------------------------------------
# 4. Column widths for the new sheet.
# ---------------------------------------------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 15.276041666666666
$ws3.Columns.Item(2).ColumnWidth = 67.16666666666667

# ---------------------------------------------------------------------------
# 5. Selection / active-tab bookkeeping: sheet3 becomes the active, selected
#    tab (mirrors the workbook.xml activeTab + sheetView tabSelected change).
# ---------------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("F10").Select()
